$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 16: shift the old F16 ("Look into other and all metrics...") content down ---
# Old state: row 15 = 2024-05-16 / 5 hrs / "added prediction metrics"; row 16 = (only) F16 TODO text.
# New state: row 15 hours become 6; a new row 16 is inserted with date 2024-05-17 / 5 hrs /
# new "confusion matrices" note, and the old F16 TODO text now lives alongside it on row 16.

# Update row 15 hours: 5 -> 6
$ws.Range("B15").Value = 6

# Preserve the existing F16 TODO text before we touch anything else (use Value2 - Value's
# getter is ambiguous on this host and returns the property descriptor instead of the data)
$todoText = $ws.Range("F16").Value2

# New row 16 data
$ws.Range("A16").Value = 45429
$ws.Range("A16").NumberFormat = $ws.Range("A15").NumberFormat
$ws.Range("B16").Value = 5
$ws.Range("D16").Value = "Looking into confusion matrices to aid in precision"
$ws.Range("F16").Value = $todoText

# --- Wrap text formatting on column D notes, matching the other wrapped note cells ---
$ws.Range("D1").WrapText = $true
$ws.Range("D2").WrapText = $true
$ws.Range("D3").WrapText = $true
$ws.Range("D4").WrapText = $true
$ws.Range("D8").WrapText = $true
$ws.Range("D9").WrapText = $true
$ws.Range("D10").WrapText = $true
$ws.Range("D11").WrapText = $true
$ws.Range("D13").WrapText = $true
$ws.Range("D15").WrapText = $true
$ws.Range("D16").WrapText = $true

# Row heights for newly-wrapped two-line rows
$ws.Rows(9).RowHeight = 28.5
$ws.Rows(11).RowHeight = 28.5

# --- View state: scroll so row 10 is at top, select D16 ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D16").Select()
